# Actualizacion Datos Personales 4 nov
# Updates the statistics row for "Saucedo Rivalcoba Graciela" (row 12)
# on both the "1er Parcial" and "3er Parcial" sheets. The "2o Parcial"
# sheet is left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("E12").Value = 28
    $ws.Range("F12").Value = 12
    $ws.Range("G12").Value = 70
    $ws.Range("H12").Value = 30
    $ws.Range("I12").Value = 8.300000000000001
    $ws.Range("J12").Value = 12
    $ws.Range("K12").Value = 30
}
